# Refactor the sadness_chat sample data:
#  - Rewrite two of the "A" messages to be more directly about sadness
#  - Remove the afternoon/evening rows (6-11), keeping only the morning
#    exchange (rows 1-5)
#  - Update the sheet selection to reflect the now-empty trailing rows
#  - Configure the page setup (A4 / portrait) for the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the wording of the two morning messages from "A"
$ws.Range("A2").Value = 'A,"오늘 아침에 창밖을 보는데 비가 내리더라… 슬퍼졌어",08:50'
$ws.Range("A4").Value = 'A,"출근길에 들었던 노래가 너무 슬프더라",09:10'

# Drop the rest of the conversation (rows 6 through 11)
$ws.Range("A6:A11").EntireRow.Delete()

# Reflect the new empty area in the sheet's selection
$ws.Range("A6:A11").Select()

# Configure page setup (A4, portrait)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
